# Add a new "6.C.2" (Waste incineration - Other, non-biogenic) mapping row
# to the "map" sheet, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Insert a new blank row at row 73 (shifts rows 73+ down to 74+).
$ws.Rows.Item(73).Insert()

# Populate the new row. Order matters for shared-string allocation:
# A (inv_sector code), then C (scaling_sector), then B (description).
$ws.Range("A73").Value = "6.C.2"
$ws.Range("C73").Value = "Waste-incineration"
$ws.Range("B73").Value = "Waste incineration - Other (non-biogenic)"
$ws.Range("D73").Value = "5C_Waste-incineration"

# Match the formatting used on similar "full record" rows (e.g. row 71),
# where the scaling_sector/ceds_sector cells carry the explicit-black-font
# style instead of the default theme color.
$ws.Range("C73").Font.Color = 0

# The edited workbook ends up with the "map" tab active/selected (rather
# than "year", which was active before).
$ws.Activate()
